$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 7.787422666666667
$ws.Range("H2").Value = 23.362268
$ws.Range("I2").Value = 0.2161047632645357
$ws.Range("J2").Value = 0.2161047632645357
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 14.94333467541022
$ws.Range("R2").Value = 134.490012078692
$ws.Range("S2").Value = 0.001409871616647264
$ws.Range("T2").Value = 0.001409871616647264
$ws.Range("G3").Value = 7.787422666666667
$ws.Range("H3").Value = 23.362268
$ws.Range("I3").Value = 0.2161047632645357
$ws.Range("J3").Value = 0.2161047632645357
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 1411.769327872737
$ws.Range("R3").Value = 12705.92395085463
$ws.Range("S3").Value = 0.1331974119468964
$ws.Range("T3").Value = 0.1331974119468964
$ws.Range("G4").Value = 7.787422666666667
$ws.Range("H4").Value = 23.362268
$ws.Range("I4").Value = 0.2161047632645357
$ws.Range("J4").Value = 0.2161047632645357
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 226.2065019425062
$ws.Range("R4").Value = 2035.858517482556
$ws.Range("S4").Value = 0.02134209890343962
$ws.Range("T4").Value = 0.02134209890343962
$ws.Range("G5").Value = 7.787422666666667
$ws.Range("H5").Value = 23.362268
$ws.Range("I5").Value = 0.2161047632645357
$ws.Range("J5").Value = 0.2161047632645357
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 637.591378654921
$ws.Range("R5").Value = 5738.322407894288
$ws.Range("S5").Value = 0.06015538079755241
$ws.Range("T5").Value = 0.06015538079755241
$ws.Range("G6").Value = 18.16892433333334
$ws.Range("H6").Value = 54.50677300000001
$ws.Range("I6").Value = 0.5041964793605992
$ws.Range("J6").Value = 0.5041964793605993
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 34.86446397308745
$ws.Range("R6").Value = 313.780175757787
$ws.Range("S6").Value = 0.003289387493018034
$ws.Range("T6").Value = 0.003289387493018035
$ws.Range("G7").Value = 18.16892433333334
$ws.Range("H7").Value = 54.50677300000001
$ws.Range("I7").Value = 0.5041964793605992
$ws.Range("J7").Value = 0.5041964793605993
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("Q7").Value = 3293.815064647057
$ws.Range("R7").Value = 29644.33558182351
$ws.Range("S7").Value = 0.3107643956989523
$ws.Range("T7").Value = 0.3107643956989524
$ws.Range("G8").Value = 18.16892433333334
$ws.Range("H8").Value = 54.50677300000001
$ws.Range("I8").Value = 0.5041964793605992
$ws.Range("J8").Value = 0.5041964793605993
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 527.7649606837936
$ws.Range("R8").Value = 4749.884646154142
$ws.Range("S8").Value = 0.04979349352012109
$ws.Range("T8").Value = 0.0497934935201211
$ws.Range("G9").Value = 18.16892433333334
$ws.Range("H9").Value = 54.50677300000001
$ws.Range("I9").Value = 0.5041964793605992
$ws.Range("J9").Value = 0.5041964793605993
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 1487.57169223043
$ws.Range("R9").Value = 13388.14523007387
$ws.Range("S9").Value = 0.1403492026485078
$ws.Range("T9").Value = 0.1403492026485078
$ws.Range("G10").Value = 6.195365666666667
$ws.Range("H10").Value = 18.586097
$ws.Range("I10").Value = 0.1719244078612872
$ws.Range("J10").Value = 0.1719244078612872
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 11.88832641508256
$ws.Range("R10").Value = 106.994937735743
$ws.Range("S10").Value = 0.00112163813139002
$ws.Range("T10").Value = 0.00112163813139002
$ws.Range("G11").Value = 6.195365666666667
$ws.Range("H11").Value = 18.586097
$ws.Range("I11").Value = 0.1719244078612872
$ws.Range("J11").Value = 0.1719244078612872
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 1123.147875431765
$ws.Range("R11").Value = 10108.33087888588
$ws.Range("S11").Value = 0.1059665961624092
$ws.Range("T11").Value = 0.1059665961624092
$ws.Range("G12").Value = 6.195365666666667
$ws.Range("H12").Value = 18.586097
$ws.Range("I12").Value = 0.1719244078612872
$ws.Range("J12").Value = 0.1719244078612872
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 179.9609518705166
$ws.Range("R12").Value = 1619.648566834649
$ws.Range("S12").Value = 0.01697893031630843
$ws.Range("T12").Value = 0.01697893031630844
$ws.Range("G13").Value = 6.195365666666667
$ws.Range("H13").Value = 18.586097
$ws.Range("I13").Value = 0.1719244078612872
$ws.Range("J13").Value = 0.1719244078612872
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 507.2424993174503
$ws.Range("R13").Value = 4565.182493857053
$ws.Range("S13").Value = 0.04785724325117949
$ws.Range("T13").Value = 0.04785724325117949
$ws.Range("G14").Value = 3.883692333333334
$ws.Range("H14").Value = 11.651077
$ws.Range("I14").Value = 0.1077743495135779
$ws.Range("J14").Value = 0.1077743495135779
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 7.452441815151444
$ws.Range("R14").Value = 67.07197633636299
$ws.Range("S14").Value = 0.0007031219214535059
$ws.Range("T14").Value = 0.0007031219214535059
$ws.Range("G15").Value = 3.883692333333334
$ws.Range("H15").Value = 11.651077
$ws.Range("I15").Value = 0.1077743495135779
$ws.Range("J15").Value = 0.1077743495135779
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 704.0683355435999
$ws.Range("R15").Value = 6336.615019892399
$ws.Range("S15").Value = 0.06642733928033058
$ws.Range("T15").Value = 0.06642733928033058
$ws.Range("G16").Value = 3.883692333333334
$ws.Range("H16").Value = 11.651077
$ws.Range("I16").Value = 0.1077743495135779
$ws.Range("J16").Value = 0.1077743495135779
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 112.8122223421455
$ws.Range("R16").Value = 1015.310001079309
$ws.Range("S16").Value = 0.01064359152397321
$ws.Range("T16").Value = 0.01064359152397321
$ws.Range("G17").Value = 3.883692333333334
$ws.Range("H17").Value = 11.651077
$ws.Range("I17").Value = 0.1077743495135779
$ws.Range("J17").Value = 0.1077743495135779
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 317.9753886585258
$ws.Range("R17").Value = 2861.778497926732
$ws.Range("S17").Value = 0.03000029678782062
$ws.Range("T17").Value = 0.03000029678782062
